$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 6 new "Document Type" rows (238-243) for XRechnung UBL/CII Invoice and
# CreditNote V3.0 (plain + Extension variants), per TICC-290 / TICC-291.
#
# Strategy: the existing rows 192-197 (XRechnung ...V2.3... entries) already
# carry the exact column formatting (styles) this table uses for a data row,
# including quirky per-column style indices that a fresh Value write would
# not reproduce (Excel falls back to each column's default style when a
# brand-new cell is written). So: clone rows 192:197 -> 238:243 (values +
# formats), fix up the handful of cells that actually differ for the new
# document types, then re-stamp the B/D/L column formatting (a full-range
# paste, and any later value overwrite, resets those three columns back to
# their plain column-default style) so every column's style exactly matches
# the template row.
# ---------------------------------------------------------------------------

$ws.Range("A192:L197").Copy() | Out-Null
$ws.Range("A238").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = $false

# Row heights: first 3 new rows are 1-line comments (30pt, like 192-194),
# last 3 are 2-line comments (45pt, like 195-197).
$ws.Rows.Item(238).RowHeight = 30
$ws.Rows.Item(239).RowHeight = 30
$ws.Rows.Item(240).RowHeight = 30
$ws.Rows.Item(241).RowHeight = 45
$ws.Rows.Item(242).RowHeight = 45
$ws.Rows.Item(243).RowHeight = 45

# Row 238: XRechnung UBL Invoice V3.0
$ws.Cells.Item(238, 1).Value = "XRechnung UBL Invoice V3.0"
$ws.Cells.Item(238, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::Invoice##urn:cen.eu:en16931:2017#compliant#urn:xeinkauf.de:kosit:xrechnung_3.0::2.1"
$ws.Cells.Item(238, 4).Value = "8.7"
$ws.Cells.Item(238, 8).Value = "TICC-291"

# Row 239: XRechnung UBL CreditNote V3.0
$ws.Cells.Item(239, 1).Value = "XRechnung UBL CreditNote V3.0"
$ws.Cells.Item(239, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:CreditNote-2::CreditNote##urn:cen.eu:en16931:2017#compliant#urn:xeinkauf.de:kosit:xrechnung_3.0::2.1"
$ws.Cells.Item(239, 4).Value = "8.7"
$ws.Cells.Item(239, 8).Value = "TICC-291"

# Row 240: XRechnung CII Invoice V3.0
$ws.Cells.Item(240, 1).Value = "XRechnung CII Invoice V3.0"
$ws.Cells.Item(240, 3).Value = "urn:un:unece:uncefact:data:standard:CrossIndustryInvoice:100::CrossIndustryInvoice##urn:cen.eu:en16931:2017#compliant#urn:xeinkauf.de:kosit:xrechnung_3.0::D16B"
$ws.Cells.Item(240, 4).Value = "8.7"
$ws.Cells.Item(240, 8).Value = "TICC-291"

# Row 241: XRechnung UBL Invoice V3.0 Extension
$ws.Cells.Item(241, 1).Value = "XRechnung UBL Invoice V3.0 Extension"
$ws.Cells.Item(241, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::Invoice##urn:cen.eu:en16931:2017#compliant#urn:xeinkauf.de:kosit:xrechnung_3.0#conformant#urn:xeinkauf. de:kosit:extension:xrechnung_3.0::2.1"
$ws.Cells.Item(241, 4).Value = "8.7"
$ws.Cells.Item(241, 8).Value = "TICC-290"

# Row 242: XRechnung UBL CreditNote V3.0 Extension
$ws.Cells.Item(242, 1).Value = "XRechnung UBL CreditNote V3.0 Extension"
$ws.Cells.Item(242, 3).Value = "urn:oasis:names:specification:ubl:schema:xsd:CreditNote-2::CreditNote##urn:cen.eu:en16931:2017#compliant#urn:xeinkauf.de:kosit:xrechnung_3.0#conformant#urn:xeinkauf. de:kosit:extension:xrechnung_3.0::2.1"
$ws.Cells.Item(242, 4).Value = "8.7"
$ws.Cells.Item(242, 8).Value = "TICC-290"

# Row 243: XRechnung CII Invoice V3.0 Extension
$ws.Cells.Item(243, 1).Value = "XRechnung CII Invoice V3.0 Extension"
$ws.Cells.Item(243, 3).Value = "urn:un:unece:uncefact:data:standard:CrossIndustryInvoice:100::CrossIndustryInvoice##urn:cen.eu:en16931:2017#compliant#urn:xeinkauf.de:kosit:xrechnung_3.0#conformant#urn:xeinkauf.de:kosit:extension:xrechnung_3.0::D16B"
$ws.Cells.Item(243, 4).Value = "8.7"
$ws.Cells.Item(243, 8).Value = "TICC-290"

# Re-stamp the B / D / L column formats (value writes above reset these
# three columns' cell style back to the column default).
$ws.Range("B192:B197").Copy() | Out-Null
$ws.Range("B238").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("D192:D197").Copy() | Out-Null
$ws.Range("D238").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("L192:L197").Copy() | Out-Null
$ws.Range("L238").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the view roughly where the author left it (bottom of the table).
$ws.Range("C241").Select()
